$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 63 (Tuesday) - update lunch-out/lunch-in/time-out entries
$ws.Range("D63").Value = 0.49861111111111112
$ws.Range("E63").Value = 0.51388888888888895
$ws.Range("F63").Value = 0.66319444444444442

# Row 64 (Wednesday) - update time-in/time-out entries
$ws.Range("C64").Value = 0.30208333333333331
$ws.Range("F64").Value = 0.61458333333333337

# New small note cell next to the week-starting date
$ws.Range("G61").Value = "  "

# Restore the selection to reflect where the user left off editing
$ws.Activate()
$ws.Range("F65").Select()
